$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume-change (E) columns with latest scraped values.
# Force each cell to Text format before assignment so numeric-looking strings
# (e.g. "1.00", "577.90") are stored as literal text, matching the source data,
# then restore the cells original style so no formatting/style drifts from the original.
$cellUpdates = @(
    @{ Cell = "D2"; Value = '61.997.24' }
    @{ Cell = "E2"; Value = '  -2.06%  ' }
    @{ Cell = "D3"; Value = '3.416.51' }
    @{ Cell = "E3"; Value = '  -1.49%  ' }
    @{ Cell = "D4"; Value = '1.00' }
    @{ Cell = "E4"; Value = '  +0.03%  ' }
    @{ Cell = "D5"; Value = '577.90' }
    @{ Cell = "E5"; Value = '  -0.57%  ' }
    @{ Cell = "D6"; Value = '152.81' }
    @{ Cell = "E6"; Value = '  +3.47%  ' }
    @{ Cell = "E7"; Value = '  +0.08%  ' }
    @{ Cell = "E9"; Value = '  +3.64%  ' }
    @{ Cell = "E10"; Value = '  -0.77%  ' }
    @{ Cell = "E11"; Value = '  +2.96%  ' }
    @{ Cell = "D12"; Value = '4.002.80' }
    @{ Cell = "E12"; Value = '  -1.47%  ' }
    @{ Cell = "E13"; Value = '  +0.67%  ' }
    @{ Cell = "D14"; Value = '28.65' }
    @{ Cell = "E14"; Value = '  -2.87%  ' }
    @{ Cell = "D15"; Value = '3.416.47' }
    @{ Cell = "E15"; Value = '  -1.47%  ' }
    @{ Cell = "E16"; Value = '  -0.63%  ' }
    @{ Cell = "D17"; Value = '62.056.16' }
    @{ Cell = "E17"; Value = '  -1.94%  ' }
    @{ Cell = "D18"; Value = '6.52' }
    @{ Cell = "E18"; Value = '  +1.92%  ' }
    @{ Cell = "E19"; Value = '  -0.09%  ' }
    @{ Cell = "E20"; Value = '  -4.13%  ' }
    @{ Cell = "D21"; Value = '381.89' }
    @{ Cell = "E21"; Value = '  -1.86%  ' }
    @{ Cell = "E22"; Value = '  +0.90%  ' }
    @{ Cell = "D23"; Value = '75.25' }
    @{ Cell = "E23"; Value = '  +1.05%  ' }
    @{ Cell = "E24"; Value = '  +0.10%  ' }
    @{ Cell = "D25"; Value = '3.560.07' }
    @{ Cell = "E25"; Value = '  -1.48%  ' }
    @{ Cell = "E26"; Value = '  -3.92%  ' }
    @{ Cell = "E27"; Value = '  -1.83%  ' }
    @{ Cell = "D28"; Value = '7.68' }
    @{ Cell = "E28"; Value = '  +0.11%  ' }
    @{ Cell = "E29"; Value = '  -0.01%  ' }
    @{ Cell = "E30"; Value = '  -3.61%  ' }
    @{ Cell = "E31"; Value = '  -1.12%  ' }
    @{ Cell = "E32"; Value = '  -0.02%  ' }
    @{ Cell = "D33"; Value = '1.33' }
    @{ Cell = "E33"; Value = '  -0.91%  ' }
    @{ Cell = "D34"; Value = '23.20' }
    @{ Cell = "E34"; Value = '  -1.18%  ' }
    @{ Cell = "D35"; Value = '5.51' }
    @{ Cell = "E35"; Value = '  +3.21%  ' }
    @{ Cell = "E36"; Value = '  -1.12%  ' }
    @{ Cell = "E37"; Value = '  -2.83%  ' }
    @{ Cell = "D38"; Value = '168.47' }
    @{ Cell = "E38"; Value = '  +0.47%  ' }
    @{ Cell = "D39"; Value = '30.91' }
    @{ Cell = "E39"; Value = '  -3.46%  ' }
    @{ Cell = "D40"; Value = '3.450.43' }
    @{ Cell = "E40"; Value = '  -1.59%  ' }
    @{ Cell = "E41"; Value = '  +2.45%  ' }
    @{ Cell = "D42"; Value = '42.68' }
    @{ Cell = "E42"; Value = '  +0.61%  ' }
    @{ Cell = "D43"; Value = '0.779' }
    @{ Cell = "E43"; Value = '  -1.87%  ' }
    @{ Cell = "E44"; Value = '  +0.51%  ' }
    @{ Cell = "E45"; Value = '  -3.93%  ' }
    @{ Cell = "E46"; Value = '  -4.18%  ' }
    @{ Cell = "D47"; Value = '2.552.23' }
    @{ Cell = "E47"; Value = '  -1.60%  ' }
    @{ Cell = "E48"; Value = '  +0.35%  ' }
    @{ Cell = "D49"; Value = '22.59' }
    @{ Cell = "E49"; Value = '  -2.04%  ' }
    @{ Cell = "E50"; Value = '  -6.40%  ' }
    @{ Cell = "E51"; Value = '  +0.10%  ' }
)

foreach ($update in $cellUpdates) {
    $range = $ws.Range($update.Cell)
    $originalStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $update.Value
    $range.Style = $originalStyle
}
